# Update existing sheet "figuras": A2 and C2 become text values instead of numbers
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("figuras")
$ws1.Range("A2").Value = "Huila"
$ws1.Range("C2").Value = "Top 5"

# Add two new sheets after "figuras": "other sheet" and "other sheet1"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "other sheet"
$ws2.Range("A2").Value = "Huila"
$ws2.Range("C2").Value = "Top 5"

$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "other sheet1"
